# Applies the changes described by the commit:
#   "Added Introduction to DevGuide and modify logic ppt"
#
# Concretely, for this single-slide deck (slide1.xml):
#   - Refresh the "last saved" date field (datetimeFigureOut) on the
#     slide master and every slide layout from 7/21/17 -> 3/25/2018.
#   - Rename three command boxes in the class diagram:
#       AddCommand   -> TaskCommand
#       ClearCommand -> EventCommand
#       FindCommand  -> HelpCommand

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Update the datetimeFigureOut placeholder text wherever it still
#    reads "7/21/17" (slide master + every custom/slide layout).
# ---------------------------------------------------------------------
function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shape = $container.Shapes.Item($i)
        if ($shape.HasTextFrame -eq -1) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "7/21/17") {
                $tr.Text = "3/25/2018"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li)
}

# ---------------------------------------------------------------------
# 2. Rename the command shapes on slide 1.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame -eq -1) {
        $tr = $shape.TextFrame.TextRange
        switch ($tr.Text) {
            "AddCommand"   { $tr.Text = "TaskCommand" }
            "ClearCommand" { $tr.Text = "EventCommand" }
            "FindCommand"  { $tr.Text = "HelpCommand" }
        }
    }
}
